$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text formatting (values such as
# "30.189.79" or "0.9977" must stay literal text, not be coerced to numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '30.189.79'
$ws.Range("E2").Value = '  -3.52%  '
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.919.43'
$ws.Range("E3").Value = '  -3.38%  '
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '0.9977'
$ws.Range("E4").Value = '  -0.77%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '245.64'
$ws.Range("E5").Value = '  -2.89%  '
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").Value = '0.7206'
$ws.Range("E6").Value = '  -8.65%  '
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '0.9972'
$ws.Range("E7").Value = '  -0.72%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.3241'
$ws.Range("E8").Value = '  -4.45%  '
$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D9").Value = '26.32'
$ws.Range("E9").Value = '  +2.19%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").Value = '43.03'
$ws.Range("E10").Value = '  -2.47%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.06836'
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '0.7889'
$ws.Range("E12").Value = '  -5.69%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07916'
$ws.Range("E13").Value = '  -2.40%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.911.21'
$ws.Range("E14").Value = '  -3.81%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '5.377'
$ws.Range("E15").Value = '  -1.21%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '94.03'
$ws.Range("E16").Value = '  -6.47%  '
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = '14.41'
$ws.Range("E17").Value = '  +3.61%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").Value = '258.86'
$ws.Range("E18").Value = '  -4.71%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '30.202.55'
$ws.Range("E19").Value = '  -3.58%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '5.811'
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.000007858'
$ws.Range("E21").Value = '  -0.89%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.165.19'
$ws.Range("E22").Value = '  -3.55%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '0.9973'
$ws.Range("E23").Value = '  -0.85%  '
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").Value = '0.9973'
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").Value = '6.837'
$ws.Range("E25").Value = '  -1.40%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '9.621'
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '159.81'
$ws.Range("E27").Value = '  -3.16%  '
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").Value = '0.1325'
$ws.Range("E28").Value = '  -9.97%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '18.71'
$ws.Range("E29").Value = '  -4.99%  '
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = '2.216'
$ws.Range("E30").Value = '  +1.99%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '1.355'
$ws.Range("E31").Value = '  -0.31%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.541'
$ws.Range("E32").Value = '  -1.60%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.384'
$ws.Range("E33").Value = '  -3.82%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.178'
$ws.Range("E34").Value = '  -3.22%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.05014'
$ws.Range("E35").Value = '  -3.13%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.187'
$ws.Range("E36").Value = '  -2.06%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.7382'
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '2.715'
$ws.Range("E38").Value = '  -3.01%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01935'
$ws.Range("E39").Value = '  -3.16%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.793'
$ws.Range("E40").Value = '  -4.46%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '79.58'
$ws.Range("E41").Value = '  +1.65%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '6.440'
$ws.Range("E42").Value = '  -2.89%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4402'
$ws.Range("E43").Value = '  -5.07%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '2.010'
$ws.Range("E44").Value = '  -2.49%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '0.9977'
$ws.Range("E45").Value = '  -0.67%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = '0.8307'
$ws.Range("E46").Value = '  -2.32%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '102.14'
$ws.Range("E47").Value = '  -3.17%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.622'
$ws.Range("E48").Value = '  -3.43%  '
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '7.225'
$ws.Range("E49").Value = '  -3.31%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '35.99'
$ws.Range("E50").Value = '  -1.30%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.4075'
$ws.Range("E51").Value = '  -4.65%  '
